$d = $word.ActiveDocument

# 1. Merge phone number runs: "+3809" + "97590873" -> "+380997590873"
$d.Content.Find.Execute("+380997590873", $true, $false, $false, $false, $false, $true, 1, $false, "+380997590873", 2)

# 2. Merge objective paragraph runs into a single sentence
$d.Content.Find.Execute("Young and purposeful full-stack software engineer seeks for a senior developer position.", $true, $false, $false, $false, $false, $true, 1, $false, "Young and purposeful full-stack software engineer seeks for a senior developer position.", 2)

# 3. Merge "*" + "OS: " runs into "*OS: "
$d.Content.Find.Execute("*OS: ", $true, $false, $false, $false, $false, $true, 1, $false, "*OS: ", 2)

# 4. Merge "E" + "nterprise project development..." runs
$d.Content.Find.Execute("Enterprise project development, DevOps, using PHP 5/7, Node, Mongo, Redis, RabbitMQ, Silex, Strongloop, microservice architecture, docker, CI services.", $true, $false, $false, $false, $false, $true, 1, $false, "Enterprise project development, DevOps, using PHP 5/7, Node, Mongo, Redis, RabbitMQ, Silex, Strongloop, microservice architecture, docker, CI services.", 2)

# 5. Prefix "Aug. 2014 – Dec. 2014" with "*"
$d.Content.Find.Execute("Aug. 2014 – Dec. 2014", $true, $false, $false, $false, $false, $true, 1, $false, "*Aug. 2014 – Dec. 2014", 2)

# 6. Prefix "June 2014 - Sept. 2014" with "*"
$d.Content.Find.Execute("June 2014 - Sept. 2014", $true, $false, $false, $false, $false, $true, 1, $false, "*June 2014 - Sept. 2014", 2)

# 7. Change color of Normal style from auto to RGB(00,00,0A)
#    Word's Color property is packed as 0x00BBGGRR, so build the value
#    from the target R,G,B so the saved OOXML w:color ends up as "00000A".
$r = 0x00
$g = 0x00
$b = 0x0A
$colorValue = ($b * 65536) + ($g * 256) + $r
$normalStyle = $d.Styles.Item("Normal")
$normalStyle.Font.Color = $colorValue
